$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# --- Shift existing block rows down to make room for the expanded layout ---
# Original rows 15-20 (Margarine block) need to become rows 18-23.
# Original rows 22-23 (Brood block) need to become rows 28-29.
$ws.Rows("15:17").Insert()
$ws.Rows("25:27").Insert()

# --- Make the two "g per leerling gebruikt" result cells bold (new style) ---
$ws.Range("A13").Font.Bold = $true
$ws.Range("A23").Font.Bold = $true

# --- New block: "Berekening voor nieuwe klas" ---
# The order in which brand-new text values are assigned below matches the
# order their entries were appended to the shared string table in the
# target workbook (A34 header, then the B36:B39 labels, then B35, then
# B40:B41) so the resulting sharedStrings.xml indices line up exactly.
$ws.Range("A34").Value = "Berekening voor nieuwe klas"
$ws.Range("A34").Font.Bold = $true
$ws.Range("A34").Font.Size = 14
$ws.Rows(34).RowHeight = 18.75

$ws.Range("A35").Value = 30

$ws.Range("A36").Formula = "=A23*A35"
$ws.Range("B36").Value = "g margarine"

$ws.Range("A37").Formula = "=ROUNDUP(A36/250,0)"
$ws.Range("B37").Value = "kuipjes van 250 gram"

$ws.Range("A38").Formula = "=A13*A35"
$ws.Range("A38").NumberFormat = "0"
$ws.Range("B38").Value = "g hagelslag"

$ws.Range("A39").Formula = "=ROUNDUP(A38/250,0)"
$ws.Range("B39").Value = "doosjes hagelslag van 250 gram"

$ws.Range("B35").Value = "leerlingen + begeleiders"

$ws.Range("A40").Formula = "=A35"
$ws.Range("B40").Value = "sneetjes"

$ws.Range("A41").Formula = "=ROUNDUP(A40/A29,0)"
$ws.Range("B41").Value = "broden"

# --- Selection moves to the last entered cell ---
$ws.Range("E41").Select()
